$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 167.5
$ws.Range("I12").Value = 127.14286
$ws.Range("K12").Value = 127.14286
$ws.Range("M12").Value = 42.85714

$ws.Range("H42").Value = 68.833336
$ws.Range("I42").Value = 53.75
$ws.Range("J42").Value = 99
$ws.Range("K42").Value = 161.25
$ws.Range("L42").Value = 297
$ws.Range("M42").Value = 68.75
$ws.Range("N42").Value = -757

$ws.Range("H74").Value = 11009.895
$ws.Range("I74").Value = 9699.3125
$ws.Range("J74").Value = 17999.666
$ws.Range("K74").Value = 9699.3125
$ws.Range("L74").Value = 17999.666
$ws.Range("M74").Value = -8763.3125
$ws.Range("N74").Value = -19871.666

$ws.Range("H76").Value = 83340120
$ws.Range("I76").Value = 6687
$ws.Range("K76").Value = 6687
$ws.Range("M76").Value = -6372

$ws.Range("H77").Value = 11009.895
$ws.Range("I77").Value = 9699.3125
$ws.Range("J77").Value = 17999.666
$ws.Range("K77").Value = 48496.5625
$ws.Range("L77").Value = 89998.33
$ws.Range("M77").Value = -43816.5625
$ws.Range("N77").Value = -99358.33

$ws.Range("H79").Value = 83340120
$ws.Range("I79").Value = 6687
$ws.Range("K79").Value = 6687
$ws.Range("M79").Value = -5595

$ws.Range("H100").Value = 2483.1667
$ws.Range("I100").Value = 1979.8
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 1979.8
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -1438.8
$ws.Range("N100").Value = -6082

$ws.Range("H138").Value = 7025.231
$ws.Range("J138").Value = 7329.4346
$ws.Range("L138").Value = 21988.3038
$ws.Range("N138").Value = -32268.3038

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 876.73914
$ws.Range("I97").Value = 876.73914
$ws.Range("K97").Value = 876.73914
$ws.Range("M97").Value = -380.73914

$ws.Range("H132").Value = 5091.8335
$ws.Range("I132").Value = 5246.316
$ws.Range("K132").Value = 15738.948
$ws.Range("M132").Value = -13208.948

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3403601.2
$ws.Range("I86").Value = 4253251.5
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 4253251.5
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -4252128.5
$ws.Range("N86").Value = -7246

$ws.Range("H89").Value = 3403601.2
$ws.Range("I89").Value = 4253251.5
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 21266257.5
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -21260641.5
$ws.Range("N89").Value = -36232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 351.46155
$ws.Range("I22").Value = 297.18182
$ws.Range("K22").Value = 297.18182
$ws.Range("M22").Value = 52.81817999999998

$ws.Range("H86").Value = 6060.5
$ws.Range("I86").Value = 4999
$ws.Range("K86").Value = 4999
$ws.Range("M86").Value = -3876

$ws.Range("H89").Value = 6060.5
$ws.Range("I89").Value = 4999
$ws.Range("K89").Value = 24995
$ws.Range("M89").Value = -19379

$ws.Range("H132").Value = 2562.3333
$ws.Range("I132").Value = 2391.923
$ws.Range("J132").Value = 3005.4
$ws.Range("K132").Value = 7175.768999999999
$ws.Range("L132").Value = 9016.200000000001
$ws.Range("M132").Value = -4645.768999999999
$ws.Range("N132").Value = -14076.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 11013
$ws.Range("I17").Value = 11013
$ws.Range("K17").Value = 33039
$ws.Range("M17").Value = -32870

$ws.Range("H19").Value = 49
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

$ws.Range("H117").Value = 864.4545000000001
$ws.Range("I117").Value = 669.8
$ws.Range("J117").Value = 1026.6666
$ws.Range("K117").Value = 2009.4
$ws.Range("L117").Value = 3079.9998
$ws.Range("M117").Value = 1432.6
$ws.Range("N117").Value = -9963.9998

$ws.Range("H120").Value = 9987.666999999999
$ws.Range("I120").Value = 9987.666999999999
$ws.Range("K120").Value = 29963.001
$ws.Range("M120").Value = -25125.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 329.82352
$ws.Range("I2").Value = 357.7143
$ws.Range("J2").Value = 199.66667
$ws.Range("K2").Value = 357.7143
$ws.Range("L2").Value = 199.66667
$ws.Range("M2").Value = -244.7143
$ws.Range("N2").Value = -425.66667

$ws.Range("H70").Value = 13877
$ws.Range("I70").Value = 11501.75
$ws.Range("J70").Value = 16252.25
$ws.Range("K70").Value = 11501.75
$ws.Range("L70").Value = 16252.25
$ws.Range("M70").Value = -11231.75
$ws.Range("N70").Value = -16792.25

$ws.Range("H73").Value = 13877
$ws.Range("I73").Value = 11501.75
$ws.Range("J73").Value = 16252.25
$ws.Range("K73").Value = 11501.75
$ws.Range("L73").Value = 16252.25
$ws.Range("M73").Value = -10565.75
$ws.Range("N73").Value = -18124.25

$ws.Range("H122").Value = 6460
$ws.Range("J122").Value = 6460
$ws.Range("L122").Value = 19380
$ws.Range("N122").Value = -24280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2956.3
$ws.Range("I40").Value = 3000.724
$ws.Range("K40").Value = 3000.724
$ws.Range("M40").Value = -2864.724

$ws.Range("H46").Value = 2231.5334
$ws.Range("J46").Value = 1831.5
$ws.Range("L46").Value = 1831.5
$ws.Range("N46").Value = -2207.5

$ws.Range("H82").Value = 3764.4
$ws.Range("I82").Value = 3825.6667
$ws.Range("J82").Value = 3672.5
$ws.Range("K82").Value = 3825.6667
$ws.Range("L82").Value = 3672.5
$ws.Range("M82").Value = -3464.6667
$ws.Range("N82").Value = -4394.5

$ws.Range("H85").Value = 3764.4
$ws.Range("I85").Value = 3825.6667
$ws.Range("J85").Value = 3672.5
$ws.Range("K85").Value = 3825.6667
$ws.Range("L85").Value = 3672.5
$ws.Range("M85").Value = -2577.6667
$ws.Range("N85").Value = -6168.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 336536.66
$ws.Range("I96").Value = 336536.66
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 336536.66
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -335163.66
$ws.Range("N96").ClearContents()

$ws.Range("H126").Value = 2125.75
$ws.Range("I126").Value = 2099.5
$ws.Range("K126").Value = 6298.5
$ws.Range("M126").Value = -3828.5

Write-Output "done"